$wb = $excel.ActiveWorkbook

# Update the Date value in the Metadata sheet
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-08-20T14:11:18+00:00"

# Remove the last data row (row 11) in the Concepts sheet
# (LIPRESCPOS / Libellé textuel de la posologie)
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Rows.Item(11).Delete()
